$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.530.41"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "2.439.65"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.12%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("D9").Value = "2.435.28"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("E10").Value = "  -3.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.156"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.350"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000175"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("D16").Value = "2.884.20"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "62.426.92"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "2.438.15"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.22%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "626.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.55%  "
$ws.Range("D28").Value = "0.0₃0960"
$ws.Range("E28").Value = "  -6.42%  "
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.137"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.375"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "146.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0524"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.598"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0229"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.00%  "
